$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D:D").Insert()

# Copy the number formatting/style from the (now shifted) column E onto the
# new column D, restricted to the rows that actually contain data, so that
# blank separator rows (5,6,37,79) don't gain a spurious D cell.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 157100
$ws.Range("D9").Value = 50800
$ws.Range("D10").Value = 106300
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 71800
$ws.Range("D15").Value = 35500
$ws.Range("D17").Value = 167400
$ws.Range("D18").Value = -10300
$ws.Range("D20").Value = 1600
$ws.Range("D21").Value = 26800
$ws.Range("D22").Value = 48700
$ws.Range("D23").Value = -57300
$ws.Range("D24").Value = 100
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -57400
$ws.Range("D27").Value = -60400
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -1600
$ws.Range("D33").Value = -60400
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -60400
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 82100
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 2700
$ws.Range("D44").Value = 5800
$ws.Range("D45").Value = 8400
$ws.Range("D46").Value = 99000
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 1112800
$ws.Range("D49").Value = 5400
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 16300
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1233500
$ws.Range("D57").Value = 12900
$ws.Range("D58").Value = 64100
$ws.Range("D59").Value = 18500
$ws.Range("D60").Value = 95500
$ws.Range("D61").Value = 813100
$ws.Range("D62").Value = 8500
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 917100
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -196100
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 316400
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -60400
$ws.Range("D83").Value = 35500
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 47700
$ws.Range("D91").Value = -11700
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 24200
$ws.Range("D96").Value = -3100
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -55200
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 16800

